# Applies the updates described in the commit:
#  - fixed timer condition in detect-update (affects DAILY_DELTA/BASE_TESTS
#    recomputation for many historical rows)
#  - changed get-data code for daily pos rate (adjusts B/C/D/E/F/G/I/J/K/M
#    for the most recent rows)
#  - removed trends data collection from build-data (adds the new
#    2020-05-16 row at the end of the sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 1729

$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 2603

$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 2689

$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 2416

$ws.Range("C42").Value = 0
$ws.Range("D42").Value = 2717

$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 1518

$ws.Range("C45").Value = 0
$ws.Range("D45").Value = 3458

$ws.Range("B47").Value = 3262
$ws.Range("C47").Value = 1
$ws.Range("D47").Value = 3261
$ws.Range("M47").Value = 62105

$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 3907
$ws.Range("M48").Value = 66012

$ws.Range("C49").Value = 0
$ws.Range("D49").Value = 4022
$ws.Range("M49").Value = 70034

$ws.Range("M50").Value = 71738

$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 1098
$ws.Range("M51").Value = 72836

$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 4287
$ws.Range("M52").Value = 77123

$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 3825
$ws.Range("M53").Value = 80948

$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 3705
$ws.Range("M54").Value = 84653

$ws.Range("B55").Value = 4239
$ws.Range("C55").Value = 1
$ws.Range("D55").Value = 4238
$ws.Range("M55").Value = 88892

$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 3491
$ws.Range("M56").Value = 92383

$ws.Range("C57").Value = 0
$ws.Range("D57").Value = 2216
$ws.Range("M57").Value = 94599

$ws.Range("C58").Value = 0
$ws.Range("D58").Value = 1785
$ws.Range("M58").Value = 96384

$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 4274
$ws.Range("M59").Value = 100658

$ws.Range("B60").Value = 4570
$ws.Range("C60").Value = 29
$ws.Range("D60").Value = 4541
$ws.Range("M60").Value = 105228

$ws.Range("B61").Value = 4637
$ws.Range("C61").Value = 1
$ws.Range("D61").Value = 4636
$ws.Range("M61").Value = 109865

$ws.Range("B62").Value = 4686
$ws.Range("C62").Value = 1
$ws.Range("D62").Value = 4685
$ws.Range("M62").Value = 114551

$ws.Range("B63").Value = 5760
$ws.Range("C63").Value = 8
$ws.Range("D63").Value = 5752
$ws.Range("M63").Value = 120311

$ws.Range("B64").Value = 2837
$ws.Range("C64").Value = 15
$ws.Range("M64").Value = 123148

$ws.Range("M65").Value = 124803

$ws.Range("B66").Value = 5212
$ws.Range("C66").Value = 3
$ws.Range("D66").Value = 5209
$ws.Range("M66").Value = 130015

$ws.Range("B67").Value = 5797
$ws.Range("C67").Value = 25
$ws.Range("D67").Value = 5772
$ws.Range("M67").Value = 135812

$ws.Range("B68").Value = 5580
$ws.Range("C68").Value = 10
$ws.Range("D68").Value = 5570
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 39
$ws.Range("M68").Value = 141392

$ws.Range("B69").Value = 5800
$ws.Range("C69").Value = 17
$ws.Range("D69").Value = 5783
$ws.Range("E69").Value = 642
$ws.Range("G69").Value = 642
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 33
$ws.Range("K69").Value = 23079
$ws.Range("M69").Value = 147192

$ws.Range("B70").Value = 5832
$ws.Range("C70").Value = 12
$ws.Range("D70").Value = 5820
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 593
$ws.Range("K70").Value = 23672
$ws.Range("M70").Value = 153024

$ws.Range("B71").Value = 3100
$ws.Range("C71").Value = 14
$ws.Range("D71").Value = 3086
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 25
$ws.Range("K71").Value = 24061
$ws.Range("M71").Value = 156124

$ws.Range("C72").Value = 0
$ws.Range("D72").Value = 1227
$ws.Range("E72").Value = 465
$ws.Range("F72").Value = 2
$ws.Range("K72").Value = 24526
$ws.Range("M72").Value = 157351

$ws.Range("B73").Value = 6395
$ws.Range("C73").Value = 156
$ws.Range("D73").Value = 6239
$ws.Range("E73").Value = 522
$ws.Range("F73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 37
$ws.Range("K73").Value = 25048
$ws.Range("M73").Value = 163746

$ws.Range("B74").Value = 5559
$ws.Range("C74").Value = 557
$ws.Range("D74").Value = 5002
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 376
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 35
$ws.Range("K74").Value = 25424
$ws.Range("M74").Value = 169305

$ws.Range("B75").Value = 4217
$ws.Range("C75").Value = 2109
$ws.Range("D75").Value = 2108
$ws.Range("F75").Value = 1
$ws.Range("G75").Value = 608
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 21
$ws.Range("K75").Value = 26033
$ws.Range("M75").Value = 173522

$ws.Range("B76").Value = 2803
$ws.Range("C76").Value = 2151
$ws.Range("D76").Value = 652
$ws.Range("E76").Value = 592
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 592
$ws.Range("H76").Value = 38
$ws.Range("I76").Value = 1
$ws.Range("J76").Value = 37
$ws.Range("K76").Value = 26625
$ws.Range("L76").Value = 1583
$ws.Range("M76").Value = 176325

$ws.Range("B77").Value = 855
$ws.Range("C77").Value = 739
$ws.Range("D77").Value = 116
$ws.Range("E77").Value = 648
$ws.Range("F77").Value = 6
$ws.Range("G77").Value = 642
$ws.Range("H77").Value = 18
$ws.Range("I77").Value = 4
$ws.Range("J77").Value = 14
$ws.Range("K77").Value = 27273
$ws.Range("L77").Value = 1601
$ws.Range("M77").Value = 177180

# New row 78 (2020-05-16) appended at the end of the trend table.
$ws.Range("A78").Value = "'2020-05-16"
$ws.Range("A78").ClearFormats()
$ws.Range("B78").Value = 63
$ws.Range("C78").Value = 63
$ws.Range("D78").Value = 0
$ws.Range("E78").Value = 505
$ws.Range("F78").Value = 502
$ws.Range("G78").Value = 3
$ws.Range("H78").Value = 6
$ws.Range("I78").Value = 6
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 27778
$ws.Range("L78").Value = 1607
$ws.Range("M78").Value = 177243
